# Weekly update: two new price records (one per "Calidad") were reported
# for Albahaca at Mercado Mayorista Lo Valledor de Santiago. Insert them as
# new rows 344-345, right after the existing row 343, pushing the previously
# recorded rows 344-367 down to 346-369.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows (existing rows 344:367 shift to 346:369).
$ws.Rows("344:345").Insert()

# New row 344 - Calidad "Primera"
$ws.Range("A344").Value = 6
$ws.Range("B344").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C344").Value = "Metropolitana"
$ws.Range("D344").Value = 44585
$ws.Range("E344").Value = 13
$ws.Range("F344").Value = 100112052
$ws.Range("G344").Value = "Albahaca"
$ws.Range("H344").Value = "Sin especificar"
$ws.Range("I344").Value = "Primera"
$ws.Range("J344").Value = 520
$ws.Range("K344").Value = 3000
$ws.Range("L344").Value = 3500
$ws.Range("M344").Value = 3240
$ws.Range("N344").Value = "`$/docena de matas"
$ws.Range("O344").Value = "Región Metropolitana"
$ws.Range("P344").Value = 540
$ws.Range("Q344").Value = 6
$ws.Range("R344").Value = "Hortaliza"

# New row 345 - Calidad "Segunda"
$ws.Range("A345").Value = 6
$ws.Range("B345").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C345").Value = "Metropolitana"
$ws.Range("D345").Value = 44585
$ws.Range("E345").Value = 13
$ws.Range("F345").Value = 100112052
$ws.Range("G345").Value = "Albahaca"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Segunda"
$ws.Range("J345").Value = 180
$ws.Range("K345").Value = 2000
$ws.Range("L345").Value = 2500
$ws.Range("M345").Value = 2222
$ws.Range("N345").Value = "`$/docena de matas"
$ws.Range("O345").Value = "Región Metropolitana"
$ws.Range("P345").Value = 370
$ws.Range("Q345").Value = 6
$ws.Range("R345").Value = "Hortaliza"
